$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the Binance conversion lines in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.61 = 9635.77 pesos`n✅ 9635.77 pesos = 2.58 = 930.27 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Sheet "tasas": update N10, O10, N12, O12 ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 383
$wsTasas.Range("O10").Value = 3690.5
$wsTasas.Range("N12").Value = 3728.9
$wsTasas.Range("O12").Value = 360
